$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 3, 6, 11, 12 per repulled data / mean recalculation
$ws.Range("F3").Value = 4
$ws.Range("F6").Value = -5
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 1
